$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Password value for the first test case (Valid Username and Valid
# Password Login Test) so the automation script can pick up the new
# credential used while getting Username/Password and handling login alerts.
$ws.Range("C2").Value = "leo_12345"

# Leave the selection where the automation's last read/write ended up.
$ws.Range("H9").Select()
